$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, derived from the source diff.
$updates = [ordered]@{
    "D2" = '60.106.75'
    "E2" = '  +0.23%  '
    "D3" = '2.416.35'
    "E3" = '  -0.16%  '
    "E4" = '  -0.05%  '
    "D5" = '554.17'
    "E5" = '  +0.33%  '
    "D6" = '136.85'
    "E6" = '  -0.33%  '
    "E7" = '  -0.01%  '
    "E8" = '  +0.37%  '
    "E9" = '  -1.08%  '
    "D10" = '5.66'
    "E10" = '  -2.35%  '
    "E11" = '  -0.70%  '
    "D12" = '0.353'
    "E12" = '  -1.21%  '
    "D13" = '24.77'
    "E13" = '  -0.27%  '
    "D14" = '2.847.23'
    "E14" = '  -0.19%  '
    "D15" = '59.989.18'
    "E15" = '  +0.12%  '
    "E16" = '  -0.16%  '
    "D17" = '2.413.64'
    "E17" = '  +0.61%  '
    "D18" = '11.23'
    "E18" = '  -0.96%  '
    "D19" = '4.52'
    "E19" = '  +3.16%  '
    "D20" = '327.42'
    "E20" = '  -1.24%  '
    "E21" = '  +1.00%  '
    "D22" = '0.999'
    "E22" = '  -0.11%  '
    "D23" = '64.75'
    "E23" = '  -1.34%  '
    "D24" = '0.180'
    "E24" = '  +5.83%  '
    "D25" = '8.64'
    "E25" = '  +0.20%  '
    "E27" = '  +5.39%  '
    "D28" = '0.0₃0773'
    "E28" = '  -1.27%  '
    "D29" = '1.77'
    "E29" = '  +0.09%  '
    "D30" = '170.82'
    "E30" = '  +0.38%  '
    "E31" = '  -1.67%  '
    "D32" = '1.07'
    "E32" = '  +4.69%  '
    "E33" = '  -3.36%  '
    "D34" = '18.49'
    "E34" = '  -0.78%  '
    "E35" = '  +0.02%  '
    "D36" = '1.32'
    "E36" = '  +2.08%  '
    "D37" = '4.23'
    "E37" = '  +1.14%  '
    "E38" = '  +0.01%  '
    "D39" = '325.08'
    "E39" = '  +3.80%  '
    "E40" = '  -1.10%  '
    "D41" = '146.45'
    "E41" = '  +4.41%  '
    "E42" = '  -1.14%  '
    "E43" = '  -0.33%  '
    "D44" = '19.81'
    "E44" = '  +2.92%  '
    "D45" = '0.0516'
    "E45" = '  -0.68%  '
    "D46" = '0.578'
    "E46" = '  +0.26%  '
    "E47" = '  -1.46%  '
    "D48" = '11.06'
    "E48" = '  +0.05%  '
    "E49" = '  -0.71%  '
    "D51" = '0.938'
    "E51" = '  -1.85%  '
}

foreach ($addr in $updates.Keys) {
    $val = $updates[$addr]
    $rng = $ws.Range($addr)
    # Force a text number format so values that look numeric (e.g. "0.180",
    # "554.17") are stored as text, matching the original inlineStr cells,
    # rather than being auto-converted into numbers by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    # Reset to the default "Normal" style so no stray style index is left
    # on the cell (keeps output identical in formatting to the source).
    $rng.Style = "Normal"
}
